$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.388106333333333
$ws.Range("H2").Value = 28.164319
$ws.Range("I2").Value = 0.2414596449149976
$ws.Range("J2").Value = 0.2414596449149975
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 723.7246422393144
$ws.Range("R2").Value = 6513.52178015383
$ws.Range("S2").Value = 0.05804247204788295
$ws.Range("T2").Value = 0.05804247204788294
$ws.Range("G3").Value = 9.388106333333333
$ws.Range("H3").Value = 28.164319
$ws.Range("I3").Value = 0.2414596449149976
$ws.Range("J3").Value = 0.2414596449149975
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 953.6441918293031
$ws.Range("R3").Value = 8582.797726463728
$ws.Range("S3").Value = 0.07648194232631231
$ws.Range("T3").Value = 0.07648194232631228
$ws.Range("G4").Value = 9.388106333333333
$ws.Range("H4").Value = 28.164319
$ws.Range("I4").Value = 0.2414596449149976
$ws.Range("J4").Value = 0.2414596449149975
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 1333.362600443266
$ws.Range("R4").Value = 12000.26340398939
$ws.Range("S4").Value = 0.1069352305408023
$ws.Range("T4").Value = 0.1069352305408023
$ws.Range("I5").Value = 0.5770971896641285
$ws.Range("J5").Value = 0.5770971896641284
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 1729.727786496233
$ws.Range("R5").Value = 15567.5500784661
$ws.Range("S5").Value = 0.1387235846875523
$ws.Range("T5").Value = 0.1387235846875523
$ws.Range("I6").Value = 0.5770971896641285
$ws.Range("J6").Value = 0.5770971896641284
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.1827945783325689
$ws.Range("T6").Value = 0.1827945783325688
$ws.Range("I7").Value = 0.5770971896641285
$ws.Range("J7").Value = 0.5770971896641284
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.2555790266440073
$ws.Range("T7").Value = 0.2555790266440072
$ws.Range("I8").Value = 0.181443165420874
$ws.Range("J8").Value = 0.1814431654208739
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 543.8378327244633
$ws.Range("R8").Value = 4894.54049452017
$ws.Range("S8").Value = 0.04361561063724712
$ws.Range("T8").Value = 0.04361561063724711
$ws.Range("I9").Value = 0.181443165420874
$ws.Range("J9").Value = 0.1814431654208739
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("S9").Value = 0.05747182191917855
$ws.Range("T9").Value = 0.05747182191917853
$ws.Range("I10").Value = 0.181443165420874
$ws.Range("J10").Value = 0.1814431654208739
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.08035573286444832
$ws.Range("T10").Value = 0.0803557328644483
